# Auto-generated Excel COM-interop script
# Applies scraped-data refresh changes (per commit message: Horarios actualizados Linea 141 - 1232)
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 09:31:15'
$ws.Cells.Item(3, 1).Value = 'Total filas: 130'
$ws.Cells.Item(28, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(29, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(49, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(50, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(59, 3).Value = '215A_EL PATO'
$ws.Cells.Item(60, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(62, 1).Value = '06:46:06'
$ws.Cells.Item(62, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(62, 4).Value = 77
$ws.Cells.Item(63, 1).Value = '08:00:50'
$ws.Cells.Item(63, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(63, 4).Value = 3
$ws.Cells.Item(89, 1).Value = '07:48:14'
$ws.Cells.Item(89, 3).Value = '215B_EL PATO'
$ws.Cells.Item(89, 4).Value = 72
$ws.Cells.Item(90, 1).Value = '08:30:59'
$ws.Cells.Item(90, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(90, 4).Value = 30
$ws.Cells.Item(96, 1).Value = '08:48:29'
$ws.Cells.Item(96, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(96, 4).Value = 26
$ws.Cells.Item(97, 1).Value = '08:00:50'
$ws.Cells.Item(97, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(97, 4).Value = 74
$ws.Cells.Item(105, 1).Value = '09:31:15'
$ws.Cells.Item(105, 2).Value = '09:31'
$ws.Cells.Item(105, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(105, 4).Value = 0
$ws.Cells.Item(107, 1).Value = '08:30:59'
$ws.Cells.Item(107, 2).Value = '09:33'
$ws.Cells.Item(107, 4).Value = 63
$ws.Cells.Item(108, 1).Value = '09:31:15'
$ws.Cells.Item(108, 2).Value = '09:34'
$ws.Cells.Item(108, 3).Value = '15_ABASTO'
$ws.Cells.Item(108, 4).Value = 3
$ws.Cells.Item(109, 1).Value = '09:31:15'
$ws.Cells.Item(109, 2).Value = '09:41'
$ws.Cells.Item(109, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(109, 4).Value = 10
$ws.Cells.Item(110, 1).Value = '08:48:29'
$ws.Cells.Item(110, 2).Value = '09:44'
$ws.Cells.Item(110, 3).Value = '14_ABASTO'
$ws.Cells.Item(110, 4).Value = 56
$ws.Cells.Item(111, 1).Value = '09:31:15'
$ws.Cells.Item(111, 2).Value = '09:45'
$ws.Cells.Item(111, 3).Value = '14_ABASTO'
$ws.Cells.Item(111, 4).Value = 14
$ws.Cells.Item(112, 1).Value = '08:30:59'
$ws.Cells.Item(112, 2).Value = '09:48'
$ws.Cells.Item(112, 3).Value = '15_ABASTO'
$ws.Cells.Item(112, 4).Value = 78
$ws.Cells.Item(113, 2).Value = '09:50'
$ws.Cells.Item(113, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(113, 4).Value = 80
$ws.Cells.Item(114, 1).Value = '09:31:15'
$ws.Cells.Item(114, 2).Value = '09:51'
$ws.Cells.Item(114, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(114, 4).Value = 20
$ws.Cells.Item(115, 2).Value = '09:55'
$ws.Cells.Item(115, 3).Value = '10_OLMOS'
$ws.Cells.Item(115, 4).Value = 85
$ws.Cells.Item(116, 1).Value = '09:31:15'
$ws.Cells.Item(116, 2).Value = '09:56'
$ws.Cells.Item(116, 3).Value = '10_OLMOS'
$ws.Cells.Item(116, 4).Value = 25
$ws.Cells.Item(117, 1).Value = '09:31:15'
$ws.Cells.Item(117, 2).Value = '10:01'
$ws.Cells.Item(117, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(117, 4).Value = 30
$ws.Cells.Item(118, 2).Value = '10:03'
$ws.Cells.Item(118, 3).Value = '215C_EL PATO'
$ws.Cells.Item(118, 4).Value = 93
$ws.Cells.Item(119, 1).Value = '09:31:15'
$ws.Cells.Item(119, 2).Value = '10:03'
$ws.Cells.Item(119, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(119, 4).Value = 32
$ws.Cells.Item(120, 1).Value = '09:31:15'
$ws.Cells.Item(120, 2).Value = '10:04'
$ws.Cells.Item(120, 3).Value = '215C_EL PATO'
$ws.Cells.Item(120, 4).Value = 33
$ws.Cells.Item(121, 1).Value = '09:31:15'
$ws.Cells.Item(121, 2).Value = '10:08'
$ws.Cells.Item(121, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(121, 4).Value = 37
$ws.Cells.Item(122, 1).Value = '09:31:15'
$ws.Cells.Item(122, 2).Value = '10:19'
$ws.Cells.Item(122, 3).Value = '17_ROMERO'
$ws.Cells.Item(122, 4).Value = 48
$ws.Cells.Item(123, 1).Value = '08:30:59'
$ws.Cells.Item(123, 2).Value = '10:19'
$ws.Cells.Item(123, 4).Value = 109
$ws.Cells.Item(124, 1).Value = '09:31:15'
$ws.Cells.Item(124, 2).Value = '10:20'
$ws.Cells.Item(124, 3).Value = '10_OLMOS'
$ws.Cells.Item(124, 4).Value = 49
$ws.Cells.Item(125, 1).Value = '08:56:14'
$ws.Cells.Item(125, 2).Value = '10:33'
$ws.Cells.Item(125, 3).Value = '14_ABASTO'
$ws.Cells.Item(125, 4).Value = 97
$ws.Cells.Item(125, 5).Value = 'LP1912'
$ws.Cells.Item(126, 1).Value = '09:31:15'
$ws.Cells.Item(126, 2).Value = '10:34'
$ws.Cells.Item(126, 3).Value = '14_ABASTO'
$ws.Cells.Item(126, 4).Value = 63
$ws.Cells.Item(126, 5).Value = 'LP1912'
$ws.Cells.Item(127, 1).Value = '09:31:15'
$ws.Cells.Item(127, 2).Value = '10:34'
$ws.Cells.Item(127, 3).Value = '15_ABASTO'
$ws.Cells.Item(127, 4).Value = 63
$ws.Cells.Item(127, 5).Value = 'LP1912'
$ws.Cells.Item(128, 1).Value = '08:48:29'
$ws.Cells.Item(128, 2).Value = '10:36'
$ws.Cells.Item(128, 3).Value = '14_ABASTO'
$ws.Cells.Item(128, 4).Value = 108
$ws.Cells.Item(128, 5).Value = 'LP1912'
$ws.Cells.Item(129, 1).Value = '09:31:15'
$ws.Cells.Item(129, 2).Value = '10:44'
$ws.Cells.Item(129, 3).Value = '10_OLMOS'
$ws.Cells.Item(129, 4).Value = 73
$ws.Cells.Item(129, 5).Value = 'LP1912'
$ws.Cells.Item(130, 1).Value = '09:31:15'
$ws.Cells.Item(130, 2).Value = '10:51'
$ws.Cells.Item(130, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(130, 4).Value = 80
$ws.Cells.Item(130, 5).Value = 'LP1912'
$ws.Cells.Item(131, 1).Value = '09:31:15'
$ws.Cells.Item(131, 2).Value = '10:57'
$ws.Cells.Item(131, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(131, 4).Value = 86
$ws.Cells.Item(131, 5).Value = 'LP1912'
$ws.Cells.Item(132, 1).Value = '09:31:15'
$ws.Cells.Item(132, 2).Value = '11:04'
$ws.Cells.Item(132, 3).Value = '17_ROMERO'
$ws.Cells.Item(132, 4).Value = 93
$ws.Cells.Item(132, 5).Value = 'LP1912'
$ws.Cells.Item(133, 1).Value = '09:31:15'
$ws.Cells.Item(133, 2).Value = '11:08'
$ws.Cells.Item(133, 3).Value = '225_C ROCA-H SUR'
$ws.Cells.Item(133, 4).Value = 97
$ws.Cells.Item(133, 5).Value = 'LP1912'
$ws.Cells.Item(134, 1).Value = '09:31:15'
$ws.Cells.Item(134, 2).Value = '11:19'
$ws.Cells.Item(134, 3).Value = '215C_EL PATO'
$ws.Cells.Item(134, 4).Value = 108
$ws.Cells.Item(134, 5).Value = 'LP1912'
$ws.Cells.Item(135, 1).Value = '09:31:15'
$ws.Cells.Item(135, 2).Value = '11:21'
$ws.Cells.Item(135, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(135, 4).Value = 110
$ws.Cells.Item(135, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 09:31:15'
$ws.Cells.Item(3, 1).Value = 'Total filas: 26'
$ws.Cells.Item(30, 1).Value = '09:31:15'
$ws.Cells.Item(30, 4).Value = 33
$ws.Cells.Item(31, 1).Value = '09:31:15'
$ws.Cells.Item(31, 2).Value = '11:19'
$ws.Cells.Item(31, 3).Value = '215C_EL PATO'
$ws.Cells.Item(31, 4).Value = 108
$ws.Cells.Item(31, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 09:31:15'
$ws.Cells.Item(3, 1).Value = 'Total filas: 24'
$ws.Cells.Item(22, 1).Value = '09:31:15'
$ws.Cells.Item(22, 4).Value = 25
$ws.Cells.Item(27, 1).Value = '09:31:15'
$ws.Cells.Item(27, 2).Value = '10:18'
$ws.Cells.Item(27, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(27, 4).Value = 47
$ws.Cells.Item(28, 1).Value = '08:30:59'
$ws.Cells.Item(28, 2).Value = '10:21'
$ws.Cells.Item(28, 4).Value = 111
$ws.Cells.Item(29, 1).Value = '09:31:15'
$ws.Cells.Item(29, 2).Value = '10:22'
$ws.Cells.Item(29, 3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(29, 4).Value = 51
$ws.Cells.Item(29, 5).Value = 'L6173'
